# TC03_Verify_BLP_Solutions_ContactUS.xlsx
# "adding wait to 3rd testcase"
#
# Insert a new row above the last data row of sheet1 ("TC03_Verify_BLP_Sol_
# Contact_etc") that contains a single WAIT keyword, shifting the existing
# "Storelocator" verification row down by one. Also restores the sheet/
# selection state left by the author after making the edit: sheet1 becomes
# the active/selected tab (instead of "Testdata"), with C8 selected there,
# while "Testdata" loses its tabSelected flag.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # TC03_Verify_BLP_Sol_Contact_etc
$ws2 = $wb.Worksheets.Item(2)   # Testdata

# --- Insert the new row 14 (old row 14 "Storelocator" shifts to row 15) ---
$ws1.Rows.Item(14).Insert()

# Match the bordered/normal formatting used throughout the table for this row.
$newRow = $ws1.Range("A14:E14")
$newRow.Borders.LineStyle = 1

# Only the Keyword column (B) is populated for a WAIT step.
$ws1.Range("B14").Value = "WAIT"

# --- Restore the active sheet / selection recorded in the saved file ---
$ws1.Activate() | Out-Null
$ws1.Range("C8").Select() | Out-Null

Write-Host "Inserted WAIT row at sheet1!A14:E14"
